# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.804.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "'2.086.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'234.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'0.0790"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("D12").Value = "'2.395.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "'14.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "'21.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "'0.774"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'5.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "'2.064.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'37.707.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'71.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "'228.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "'169.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'0.140"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.67%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "'19.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "'4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'0.0634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "'3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'5.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").Value = "'98.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "'1.464.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "'1.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "'15.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "'7.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").Value = "'3.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'2.280.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
